$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H2 value
$ws.Range("H2").Value = 1.39

# Delete row 4 entirely (shifts cells up, shrinks merges automatically)
$ws.Rows.Item(4).Delete()
